# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 23, pushing the existing
# rows 23..55 down to 24..56 (matches Excel's "insert entire row" behaviour,
# which also updates the used-range dimension automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new weekly record. It mirrors
# the record that is now on row 24 (the old row 23), except for the date
# (Fecha) and the volume (Volumen), which carry the new week's figures.
$ws.Cells.Item(23, 1).Value  = 10
$ws.Cells.Item(23, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(23, 3).Value  = "La Araucanía"
$ws.Cells.Item(23, 4).Value  = 44658
$ws.Cells.Item(23, 5).Value  = 9
$ws.Cells.Item(23, 6).Value  = "Fruta"
$ws.Cells.Item(23, 7).Value  = 100107
$ws.Cells.Item(23, 8).Value  = "Otros"
$ws.Cells.Item(23, 9).Value  = 100107011
$ws.Cells.Item(23, 10).Value = "Tuna"
$ws.Cells.Item(23, 11).Value = "Sin especificar"
$ws.Cells.Item(23, 12).Value = "Primera"
$ws.Cells.Item(23, 13).Value = 200
$ws.Cells.Item(23, 14).Value = 17000
$ws.Cells.Item(23, 15).Value = 17000
$ws.Cells.Item(23, 16).Value = 17000
$ws.Cells.Item(23, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(23, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(23, 19).Value = 1062
$ws.Cells.Item(23, 20).Value = 16
